$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 15300
$ws.Range("E8").Value = 3100
$ws.Range("F8").Value = 2200
$ws.Range("H8").Value = 600
$ws.Range("D12").Value = 38300
$ws.Range("E12").Value = 21300
$ws.Range("F12").Value = 18300
$ws.Range("G12").Value = 13900
$ws.Range("H12").Value = 6800
$ws.Range("F14").Value = 5000
$ws.Range("G14").Value = 5100
$ws.Range("D17").Value = 64200
$ws.Range("E17").Value = 34100
$ws.Range("F17").Value = 28100
$ws.Range("G17").Value = 21000
$ws.Range("H17").Value = 11700
$ws.Range("D18").Value = -48900
$ws.Range("E18").Value = -31000
$ws.Range("F18").Value = -25800
$ws.Range("G18").Value = -19500
$ws.Range("H18").Value = -11100
$ws.Range("D20").Value = -32800
$ws.Range("E20").Value = -21900
$ws.Range("D21").Value = -81400
$ws.Range("E21").Value = -52700
$ws.Range("F21").Value = -25800
$ws.Range("G21").Value = -19200
$ws.Range("H21").Value = -10800
$ws.Range("D23").Value = -81700
$ws.Range("E23").Value = -53000
$ws.Range("F23").Value = -26000
$ws.Range("G23").Value = -19500
$ws.Range("H23").Value = -11100
$ws.Range("D26").Value = -82000
$ws.Range("E26").Value = -53000
$ws.Range("F26").Value = -26000
$ws.Range("G26").Value = -19500
$ws.Range("H26").Value = -11100
$ws.Range("D27").Value = -82000
$ws.Range("E27").Value = -53000
$ws.Range("F27").Value = -26000
$ws.Range("G27").Value = -19500
$ws.Range("H27").Value = -11100
$ws.Range("D32").Value = 32800
$ws.Range("E32").Value = 21900
$ws.Range("D33").Value = -82000
$ws.Range("E33").Value = -53000
$ws.Range("F33").Value = -26000
$ws.Range("G33").Value = -19500
$ws.Range("H33").Value = -11100
$ws.Range("D35").Value = -82000
$ws.Range("E35").Value = -53000
$ws.Range("F35").Value = -26000
$ws.Range("G35").Value = -19500
$ws.Range("H35").Value = -11100
$ws.Range("D41").Value = 167900
$ws.Range("E41").Value = 63900
$ws.Range("F41").Value = 36900
$ws.Range("H41").Value = 11900
$ws.Range("D42").Value = 38200
$ws.Range("E42").Value = 13300
$ws.Range("D43").Value = 4200
$ws.Range("E43").Value = 1300
$ws.Range("G43").Value = 600
$ws.Range("H43").Value = 800
$ws.Range("D45").Value = 2400
$ws.Range("E45").Value = 1200
$ws.Range("D46").Value = 211100
$ws.Range("E46").Value = 79700
$ws.Range("F46").Value = 38700
$ws.Range("G46").Value = 2700
$ws.Range("H46").Value = 12900
$ws.Range("D47").Value = 7900
$ws.Range("D48").Value = 2600
$ws.Range("H48").Value = 400
$ws.Range("E49").Value = 800
$ws.Range("H49").Value = 600
$ws.Range("D52").Value = 100
$ws.Range("F52").Value = 200
$ws.Range("D54").Value = 220800
$ws.Range("E54").Value = 81100
$ws.Range("F54").Value = 39800
$ws.Range("G54").Value = 4000
$ws.Range("H54").Value = 13900
$ws.Range("D57").Value = 3200
$ws.Range("E57").Value = 2600
$ws.Range("F57").Value = 2700
$ws.Range("G57").Value = 2700
$ws.Range("D59").Value = 27800
$ws.Range("E59").Value = 5900
$ws.Range("F59").Value = 4200
$ws.Range("G59").Value = 3700
$ws.Range("H59").Value = 3200
$ws.Range("D60").Value = 31000
$ws.Range("E60").Value = 8700
$ws.Range("F60").Value = 7100
$ws.Range("G60").Value = 6500
$ws.Range("H60").Value = 4400
$ws.Range("F61").Value = 500
$ws.Range("G61").Value = 700
$ws.Range("H61").Value = 900
$ws.Range("D62").Value = 126300
$ws.Range("E62").Value = 33900
$ws.Range("F62").Value = 400
$ws.Range("D66").Value = 157300
$ws.Range("E66").Value = 42900
$ws.Range("F66").Value = 8100
$ws.Range("G66").Value = 8000
$ws.Range("H66").Value = 5300
$ws.Range("D72").Value = -178100
$ws.Range("E72").Value = -120400
$ws.Range("F72").Value = -71100
$ws.Range("G72").Value = -45700
$ws.Range("H72").Value = -26400
$ws.Range("D76").Value = 63500
$ws.Range("E76").Value = 38200
$ws.Range("F76").Value = 31800
$ws.Range("G76").Value = -4000
$ws.Range("H76").Value = 8600
$ws.Range("D81").Value = -82000
$ws.Range("E81").Value = -53000
$ws.Range("F81").Value = -26000
$ws.Range("G81").Value = -19500
$ws.Range("H81").Value = -11100
$ws.Range("D89").Value = -42000
$ws.Range("E89").Value = -28900
$ws.Range("F89").Value = -25800
$ws.Range("G89").Value = -16400
$ws.Range("H89").Value = -7800
$ws.Range("D94").Value = -46700
$ws.Range("D100").Value = 208900
$ws.Range("E100").Value = 56300
$ws.Range("F100").Value = 61000
$ws.Range("G100").Value = 6800
$ws.Range("H100").Value = 13300
$ws.Range("D101").Value = -16200
$ws.Range("D102").Value = 104100
$ws.Range("E102").Value = 27000
$ws.Range("F102").Value = 35100
$ws.Range("G102").Value = -9700
$ws.Range("H102").Value = 5500
